# "Config capture to work with maven sure-fire plugin folder"
#
# The recorded test-data workbook gets a new "alertmessage" column (with its
# expected value) on the AddCustomerTest sheet, and the cursor/window
# bookkeeping that Excel stamps on save moves along with it.

$wb = $excel.ActiveWorkbook

# The edited sheet is the first one (tab "AddCustomerTest").
$ws = $wb.Worksheets.Item("AddCustomerTest")
$ws.Activate()

# New header + expected value for the customer-added alert message.
$ws.Cells.Item(1, 4).Value = "alertmessage"
$ws.Cells.Item(2, 4).Value = "Customer added successfully"

# Excel also records the last selection and the window scroll position when
# the file is saved after this editing session.
$ws.Range("G10").Select()
$excel.ActiveWindow.Left = 5580
